# Generate Report for Handoff
#
# Status moves from "In Translation" to "Ready for handoff" and the
# handoff timestamps are refreshed. Touches the Overview sheet (status
# shown per-language) plus each language sheet's own Status / Latest
# Handoff Datetime cells. Widening columns E/F (Overview) and C (zh-cn,
# de-de) to fit the new, longer status text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Status text: "In Translation" -> "Ready for handoff" ---
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$zhcn.Range("C2").Value = "Ready for handoff"
$dede.Range("C2").Value = "Ready for handoff"

# --- Refresh handoff timestamps ---
$overview.Range("G2").Value = "2016-08-12 03:00:58"
$zhcn.Range("H2").Value = "2016-08-12 03:00:53"
$dede.Range("H2").Value = "2016-08-12 03:00:58"

# --- Widen the status columns to fit "Ready for handoff" ---
$overview.Range("E1:F1").ColumnWidth = 16.3
$zhcn.Range("C1").ColumnWidth = 16.3
$dede.Range("C1").ColumnWidth = 16.3
